$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 0.20014
$ws.Cells.Item(2, 8).Value = 0.6004200000000001
$ws.Cells.Item(2, 9).Value = 0.09409628186320101
$ws.Cells.Item(2, 10).Value = 0.09409628186320104
$ws.Cells.Item(2, 13).Value = 3.257987
$ws.Cells.Item(2, 14).Value = 9.773961
$ws.Cells.Item(2, 15).Value = 0.06563198179863573
$ws.Cells.Item(2, 16).Value = 0.06563198179863572
$ws.Cells.Item(2, 17).Value = 0.65205351818
$ws.Cells.Item(2, 18).Value = 5.868481663620001
$ws.Cells.Item(2, 19).Value = 0.006175725458564906
$ws.Cells.Item(2, 20).Value = 0.006175725458564907
$ws.Cells.Item(3, 7).Value = 0.20014
$ws.Cells.Item(3, 8).Value = 0.6004200000000001
$ws.Cells.Item(3, 9).Value = 0.09409628186320101
$ws.Cells.Item(3, 10).Value = 0.09409628186320104
$ws.Cells.Item(3, 15).Value = 0.793654179920784
$ws.Cells.Item(3, 16).Value = 0.793654179920784
$ws.Cells.Item(3, 17).Value = 7.8849516052
$ws.Cells.Item(3, 18).Value = 70.96456444680001
$ws.Cells.Item(3, 19).Value = 0.07467990741573374
$ws.Cells.Item(3, 20).Value = 0.07467990741573377
$ws.Cells.Item(4, 7).Value = 0.20014
$ws.Cells.Item(4, 8).Value = 0.6004200000000001
$ws.Cells.Item(4, 9).Value = 0.09409628186320101
$ws.Cells.Item(4, 10).Value = 0.09409628186320104
$ws.Cells.Item(4, 13).Value = 1.833856
$ws.Cells.Item(4, 14).Value = 5.501568000000001
$ws.Cells.Item(4, 15).Value = 0.03694293550383072
$ws.Cells.Item(4, 16).Value = 0.03694293550383072
$ws.Cells.Item(4, 17).Value = 0.36702793984
$ws.Cells.Item(4, 18).Value = 3.303251458560001
$ws.Cells.Item(4, 19).Value = 0.003476192872022511
$ws.Cells.Item(4, 20).Value = 0.003476192872022512
$ws.Cells.Item(5, 7).Value = 0.20014
$ws.Cells.Item(5, 8).Value = 0.6004200000000001
$ws.Cells.Item(5, 9).Value = 0.09409628186320101
$ws.Cells.Item(5, 10).Value = 0.09409628186320104
$ws.Cells.Item(5, 13).Value = 5.151212
$ws.Cells.Item(5, 14).Value = 15.453636
$ws.Cells.Item(5, 15).Value = 0.1037709027767496
$ws.Cells.Item(5, 16).Value = 0.1037709027767495
$ws.Cells.Item(5, 17).Value = 1.03096356968
$ws.Cells.Item(5, 18).Value = 9.27867212712
$ws.Cells.Item(5, 19).Value = 0.009764456116879854
$ws.Cells.Item(5, 20).Value = 0.009764456116879856
$ws.Cells.Item(6, 9).Value = 0.7283659026117116
$ws.Cells.Item(6, 10).Value = 0.7283659026117117
$ws.Cells.Item(6, 13).Value = 3.257987
$ws.Cells.Item(6, 14).Value = 9.773961
$ws.Cells.Item(6, 15).Value = 0.06563198179863573
$ws.Cells.Item(6, 16).Value = 0.06563198179863572
$ws.Cells.Item(6, 17).Value = 5.047314728235333
$ws.Cells.Item(6, 18).Value = 45.425832554118
$ws.Cells.Item(6, 19).Value = 0.04780409766295875
$ws.Cells.Item(6, 20).Value = 0.04780409766295874
$ws.Cells.Item(7, 9).Value = 0.7283659026117116
$ws.Cells.Item(7, 10).Value = 0.7283659026117117
$ws.Cells.Item(7, 15).Value = 0.793654179920784
$ws.Cells.Item(7, 16).Value = 0.793654179920784
$ws.Cells.Item(7, 19).Value = 0.5780706431195597
$ws.Cells.Item(7, 20).Value = 0.5780706431195597
$ws.Cells.Item(8, 9).Value = 0.7283659026117116
$ws.Cells.Item(8, 10).Value = 0.7283659026117117
$ws.Cells.Item(8, 13).Value = 1.833856
$ws.Cells.Item(8, 14).Value = 5.501568000000001
$ws.Cells.Item(8, 15).Value = 0.03694293550383072
$ws.Cells.Item(8, 16).Value = 0.03694293550383072
$ws.Cells.Item(8, 17).Value = 2.841032944042667
$ws.Cells.Item(8, 18).Value = 25.56929649638401
$ws.Cells.Item(8, 19).Value = 0.02690797456337391
$ws.Cells.Item(8, 20).Value = 0.02690797456337391
$ws.Cells.Item(9, 9).Value = 0.7283659026117116
$ws.Cells.Item(9, 10).Value = 0.7283659026117117
$ws.Cells.Item(9, 13).Value = 5.151212
$ws.Cells.Item(9, 14).Value = 15.453636
$ws.Cells.Item(9, 15).Value = 0.1037709027767496
$ws.Cells.Item(9, 16).Value = 0.1037709027767495
$ws.Cells.Item(9, 17).Value = 7.980322879085334
$ws.Cells.Item(9, 18).Value = 71.822905911768
$ws.Cells.Item(9, 19).Value = 0.07558318726581936
$ws.Cells.Item(9, 20).Value = 0.07558318726581936
$ws.Cells.Item(10, 7).Value = 0.25539
$ws.Cells.Item(10, 8).Value = 0.76617
$ws.Cells.Item(10, 9).Value = 0.1200721965876032
$ws.Cells.Item(10, 10).Value = 0.1200721965876032
$ws.Cells.Item(10, 13).Value = 3.257987
$ws.Cells.Item(10, 14).Value = 9.773961
$ws.Cells.Item(10, 15).Value = 0.06563198179863573
$ws.Cells.Item(10, 16).Value = 0.06563198179863572
$ws.Cells.Item(10, 17).Value = 0.83205729993
$ws.Cells.Item(10, 18).Value = 7.48851569937
$ws.Cells.Item(10, 19).Value = 0.007880576220959785
$ws.Cells.Item(10, 20).Value = 0.007880576220959785
$ws.Cells.Item(11, 7).Value = 0.25539
$ws.Cells.Item(11, 8).Value = 0.76617
$ws.Cells.Item(11, 9).Value = 0.1200721965876032
$ws.Cells.Item(11, 10).Value = 0.1200721965876032
$ws.Cells.Item(11, 15).Value = 0.793654179920784
$ws.Cells.Item(11, 16).Value = 0.793654179920784
$ws.Cells.Item(11, 17).Value = 10.0616458002
$ws.Cells.Item(11, 18).Value = 90.5548122018
$ws.Cells.Item(11, 19).Value = 0.0952958007140214
$ws.Cells.Item(11, 20).Value = 0.09529580071402141
$ws.Cells.Item(12, 7).Value = 0.25539
$ws.Cells.Item(12, 8).Value = 0.76617
$ws.Cells.Item(12, 9).Value = 0.1200721965876032
$ws.Cells.Item(12, 10).Value = 0.1200721965876032
$ws.Cells.Item(12, 13).Value = 1.833856
$ws.Cells.Item(12, 14).Value = 5.501568000000001
$ws.Cells.Item(12, 15).Value = 0.03694293550383072
$ws.Cells.Item(12, 16).Value = 0.03694293550383072
$ws.Cells.Item(12, 17).Value = 0.4683484838400001
$ws.Cells.Item(12, 18).Value = 4.21513635456
$ws.Cells.Item(12, 19).Value = 0.004435819414339109
$ws.Cells.Item(12, 20).Value = 0.004435819414339109
$ws.Cells.Item(13, 7).Value = 0.25539
$ws.Cells.Item(13, 8).Value = 0.76617
$ws.Cells.Item(13, 9).Value = 0.1200721965876032
$ws.Cells.Item(13, 10).Value = 0.1200721965876032
$ws.Cells.Item(13, 13).Value = 5.151212
$ws.Cells.Item(13, 14).Value = 15.453636
$ws.Cells.Item(13, 15).Value = 0.1037709027767496
$ws.Cells.Item(13, 16).Value = 0.1037709027767495
$ws.Cells.Item(13, 17).Value = 1.31556803268
$ws.Cells.Item(13, 18).Value = 11.84011229412
$ws.Cells.Item(13, 19).Value = 0.01246000023828293
$ws.Cells.Item(13, 20).Value = 0.01246000023828293
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.1222276666666667
$ws.Cells.Item(14, 8).Value = 0.366683
$ws.Cells.Item(14, 9).Value = 0.05746561893748399
$ws.Cells.Item(14, 10).Value = 0.057465618937484
$ws.Cells.Item(14, 13).Value = 3.257987
$ws.Cells.Item(14, 14).Value = 9.773961
$ws.Cells.Item(14, 15).Value = 0.06563198179863573
$ws.Cells.Item(14, 16).Value = 0.06563198179863572
$ws.Cells.Item(14, 17).Value = 0.3982161490403333
$ws.Cells.Item(14, 18).Value = 3.583945341363
$ws.Cells.Item(14, 19).Value = 0.003771582456152286
$ws.Cells.Item(14, 20).Value = 0.003771582456152286
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.1222276666666667
$ws.Cells.Item(15, 8).Value = 0.366683
$ws.Cells.Item(15, 9).Value = 0.05746561893748399
$ws.Cells.Item(15, 10).Value = 0.057465618937484
$ws.Cells.Item(15, 15).Value = 0.793654179920784
$ws.Cells.Item(15, 16).Value = 0.793654179920784
$ws.Cells.Item(15, 17).Value = 4.815425384646667
$ws.Cells.Item(15, 18).Value = 43.33882846182
$ws.Cells.Item(15, 19).Value = 0.04560782867146913
$ws.Cells.Item(15, 20).Value = 0.04560782867146914
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.1222276666666667
$ws.Cells.Item(16, 8).Value = 0.366683
$ws.Cells.Item(16, 9).Value = 0.05746561893748399
$ws.Cells.Item(16, 10).Value = 0.057465618937484
$ws.Cells.Item(16, 13).Value = 1.833856
$ws.Cells.Item(16, 14).Value = 5.501568000000001
$ws.Cells.Item(16, 15).Value = 0.03694293550383072
$ws.Cells.Item(16, 16).Value = 0.03694293550383072
$ws.Cells.Item(16, 17).Value = 0.2241479398826667
$ws.Cells.Item(16, 18).Value = 2.017331458944
$ws.Cells.Item(16, 19).Value = 0.002122948654095184
$ws.Cells.Item(16, 20).Value = 0.002122948654095185
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.1222276666666667
$ws.Cells.Item(17, 8).Value = 0.366683
$ws.Cells.Item(17, 9).Value = 0.05746561893748399
$ws.Cells.Item(17, 10).Value = 0.057465618937484
$ws.Cells.Item(17, 13).Value = 5.151212
$ws.Cells.Item(17, 14).Value = 15.453636
$ws.Cells.Item(17, 15).Value = 0.1037709027767496
$ws.Cells.Item(17, 16).Value = 0.1037709027767495
$ws.Cells.Item(17, 17).Value = 0.6296206232653333
$ws.Cells.Item(17, 18).Value = 5.666585609387999
$ws.Cells.Item(17, 19).Value = 0.005963259155767389
$ws.Cells.Item(17, 20).Value = 0.00596325915576739
